$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.361.35"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "3.590.63"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.68%  "
$ws.Range("D7").Value = "3.590.98"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "4.196.28"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000206"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "3.591.56"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "66.382.36"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.610"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000121"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "3.585.53"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.157"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "174.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.882"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  -5.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.950"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.25%  "
